# Fix APD experiments + better matrix completion experiments
#
# This updates the refreshed experiment results (rows 2-7) in the APD /
# matrix-completion results sheet: several columns of rerun numbers
# (fval_ANCF, t_APD, iter_APD, timings, etc.) changed because the
# experiments were redone, rows 5-7 were regenerated wholesale (new m/M
# sizes, new result columns, and "relative" opt_type), B6 picked up an
# integer display format, and the active selection / a couple of column
# widths shifted as a side effect of editing the sheet in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: partial refresh of fval_ANCF, t_ANCF, iter_UPF/APD, timings ---
$ws.Range("D2").Value = -8.501886355030367
$ws.Range("F2").Value = -8.501822050453278
$ws.Range("H2").Value = 5214
$ws.Range("J2").Value = 1081
$ws.Range("K2").Value = 117.0451698
$ws.Range("L2").Value = 34.4810608
$ws.Range("M2").Value = 18.1343358
$ws.Range("N2").Value = 4.3878938

# --- Row 3: partial refresh ---
$ws.Range("D3").Value = -7.605965997262343
$ws.Range("F3").Value = -7.604082460342722
$ws.Range("H3").Value = 11081
$ws.Range("J3").Value = 3343
$ws.Range("K3").Value = 339.6584682
$ws.Range("L3").Value = 74.6775758
$ws.Range("M3").Value = 19.7345254
$ws.Range("N3").Value = 14.3181497

# --- Row 4: partial refresh ---
$ws.Range("D4").Value = -7.52091268779711
$ws.Range("F4").Value = -7.511552514323165
$ws.Range("H4").Value = 14821
$ws.Range("J4").Value = 7110
$ws.Range("K4").Value = 557.1190712
$ws.Range("L4").Value = 98.619882
$ws.Range("M4").Value = 60.7317447
$ws.Range("N4").Value = 29.1987912

# --- Row 5: full refresh with new m/M experiment and opt_type ---
$ws.Range("A5").Value = 999.9999999841675
$ws.Range("B5").Value = 10000000
$ws.Range("C5").Value = -75.11860661045505
$ws.Range("D5").Value = -75.21101623935986
$ws.Range("E5").Value = -75.11934343753752
$ws.Range("F5").Value = -75.11585252239679
$ws.Range("G5").Value = 121618
$ws.Range("H5").Value = 14825
$ws.Range("I5").Value = 22096
$ws.Range("J5").Value = 10142
$ws.Range("K5").Value = 576.8536942
$ws.Range("L5").Value = 106.7051043
$ws.Range("M5").Value = 63.2155754
$ws.Range("N5").Value = 43.8542076
$ws.Range("P5").Value = "relative"

# --- Row 6: full refresh with new m/M experiment and opt_type ---
$ws.Range("A6").Value = 99.99999999917024
$ws.Range("B6").Value = 9999999.999091618
$ws.Range("C6").Value = -7.498860107683855
$ws.Range("D6").Value = -5.5992913022744695
$ws.Range("E6").Value = -7.502055248090166
$ws.Range("F6").Value = -7.458520253104405
$ws.Range("G6").Value = 130302
$ws.Range("H6").Value = 16582
$ws.Range("I6").Value = 32743
$ws.Range("J6").Value = 11893
$ws.Range("K6").Value = 611.5222934
$ws.Range("L6").Value = 114.7723979
$ws.Range("M6").Value = 92.4251362
$ws.Range("N6").Value = 50.775841
$ws.Range("P6").Value = "relative"
# B6's sample-size value (M) is effectively 1e7 (9999999.9990916178), so
# display it rounded to an integer.
$ws.Range("B6").NumberFormat = "0"

# --- Row 7: full refresh with new m/M experiment and opt_type ---
$ws.Range("A7").Value = 10.000000000698492
$ws.Range("B7").Value = 9999999.999999996
$ws.Range("C7").Value = -0.7487787592463879
$ws.Range("D7").Value = 0.6972425620768496
$ws.Range("E7").Value = -0.710907289290244
$ws.Range("F7").Value = -0.6810649253910285
$ws.Range("G7").Value = 230428
$ws.Range("H7").Value = 35410
$ws.Range("I7").Value = 26509
$ws.Range("J7").Value = 20375
$ws.Range("K7").Value = 1074.6186835
$ws.Range("L7").Value = 241.800633
$ws.Range("M7").Value = 73.4349055
$ws.Range("N7").Value = 85.3975312
$ws.Range("P7").Value = "relative"

# --- Column width touch-ups (m, M and the two timing columns got wider
#     to fit the new, larger numbers) ---
$ws.Columns("A").ColumnWidth = 4.333333333333333
$ws.Columns("B").ColumnWidth = 10.833333333333334
$ws.Columns("M").ColumnWidth = 9.833333333333334
$ws.Columns("N").ColumnWidth = 9.833333333333334

# --- Active cell ended up on K9 after the edits ---
$ws.Range("K9").Select()
